$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.090.09'
$ws.Range("E2").Value = '  +1.66%  '

$ws.Range("D3").Value = '1.791.07'
$ws.Range("E3").Value = '  +1.98%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.84'
$ws.Range("E5").Value = '  +1.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4526'
$ws.Range("E7").Value = '  +0.92%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3608'
$ws.Range("E8").Value = '  +0.47%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07539'
$ws.Range("E9").Value = '  +0.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.49'
$ws.Range("E10").Value = '  +1.17%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.113'
$ws.Range("E11").Value = '  +1.27%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.0000'
$ws.Range("E12").Value = '  +0.21%  '

$ws.Range("E13").Value = '  +0.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.102'
$ws.Range("E14").Value = '  +1.07%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.275'
$ws.Range("E15").Value = '  +2.24%  '

$ws.Range("D16").Value = '1.783.58'
$ws.Range("E16").Value = '  +1.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '94.40'
$ws.Range("E17").Value = '  +1.15%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001069'
$ws.Range("E18").Value = '  +0.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06464'
$ws.Range("E19").Value = '  +0.81%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9991'
$ws.Range("E20").Value = '  +0.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.26'
$ws.Range("E21").Value = '  +2.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.860'
$ws.Range("E22").Value = '  +0.53%  '

$ws.Range("D23").Value = '28.127.72'
$ws.Range("E23").Value = '  +1.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.43'
$ws.Range("E24").Value = '  +1.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.094'
$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.43'
$ws.Range("E26").Value = '  +0.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.45'
$ws.Range("E27").Value = '  -0.20%  '

$ws.Range("D28").Value = '1.991.04'
$ws.Range("E28").Value = '  +2.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.266'
$ws.Range("E29").Value = '  +8.59%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.70'
$ws.Range("E30").Value = '  +0.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.112'
$ws.Range("E31").Value = '  +2.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09215'
$ws.Range("E32").Value = '  +1.40%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.631'
$ws.Range("E33").Value = '  +1.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.675'
$ws.Range("E34").Value = '  +0.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.00'
$ws.Range("E35").Value = '  -0.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02313'
$ws.Range("E36").Value = '  +0.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06158'
$ws.Range("E37").Value = '  +2.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2107'
$ws.Range("E38").Value = '  +0.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6395'
$ws.Range("E39").Value = '  +0.20%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.013'
$ws.Range("E40").Value = '  +0.95%  '

$ws.Range("E41").Value = '  -0.65%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.398'
$ws.Range("E42").Value = '  +1.32%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.010'
$ws.Range("E43").Value = '  +2.43%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.33'
$ws.Range("E44").Value = '  +0.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5950'
$ws.Range("E45").Value = '  +0.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.747'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '123.66'
$ws.Range("E47").Value = '  +0.46%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.978'
$ws.Range("E48").Value = '  +1.01%  '

$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.151'
$ws.Range("E49").Value = '  +0.79%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06968'
$ws.Range("E50").Value = '  +1.48%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.15'
$ws.Range("E51").Value = '  +0.86%  '
